# Added dynamic messaging system
# Replace four static notification messages with templated/dynamic versions
# that reference placeholders (gender / clothing) so the same message can be
# reused with different substituted values at runtime.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Row 2  (Id=0): the stolen-phone message now varies by the suspect's gender
$ws.Range("D2").Value = "Een [geslacht] heeft hier een mobiel uit mijn winkel gestolen"

# Row 7  (Id=5): the "ran away" message now varies by the suspect's gender
$ws.Range("D7").Value = "De [geslacht] rende hier gauw weg!"

# Row 10 (Id=8): the "saw someone running" message now references a templated top garment
$ws.Range("D10").Value = "Ik zag iemand hier wegrennen met een [bovenstuk] bovenstuk!"

# Row 12 (Id=10): the "is he still around" message now references a templated bottom garment
$ws.Range("D12").Value = "Loopt die persoon met het [onderstuk] onderstuk hier nog ergens rond?"

# Column D needs to be a bit wider to comfortably fit the new, longer messages
$ws.Columns("D").ColumnWidth = 58.5

# Reflect the last-used selection in the sheet (matches the edited state)
$ws.Range("D14").Select() | Out-Null
